$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("students")

# --- Update "Masters Thesis" (M) column checkmarks ---
# Row 7 (Cole Brauer) gains a Masters Thesis "x"
$ws.Range("M7").Value = "x"
# Row 9 (Sudhanshu Katarey) loses the Masters Thesis "x"
$ws.Range("M9").Clear()

# --- Update "former" (O) column checkmarks ---
# Row 17 (Mannat Rana) and Row 20 (Clint Ewell) are now marked as former
$ws.Range("O17").Value = "x"
$ws.Range("O20").Value = "x"

# --- Update "start" (P) column ---
# Row 17 (Mannat Rana) start term changes from Spring 2019 to Fall 2019
$ws.Range("P17").Value = "Fall 2019"

# --- Populate "stop" (Q) column ---
$ws.Range("Q2").Value = "current"
$ws.Range("Q3").Value = "current"
$ws.Range("Q4").Value = "current"
$ws.Range("Q5").Value = "current"
$ws.Range("Q7").Value = "current"
$ws.Range("Q8").Value = "current"
$ws.Range("Q9").Value = "current"
$ws.Range("Q16").Value = "current"
$ws.Range("Q18").Value = "Spring 2020"
$ws.Range("Q19").Value = "Spring 2020"
$ws.Range("Q21").Value = "Spring 2020"
$ws.Range("Q22").Value = "current"
$ws.Range("Q23").Value = "current"

# --- Column M (Masters Thesis) becomes its own best-fit-width column, split out
#     from the shared default-width block (3-16384 -> 3-12 / 13 / 14-16384) ---
$ws.Columns.Item(13).ColumnWidth = 13.25

# --- Update selection to match the new active cell ---
$ws.Range("Q23").Select()
